# Insert a new record row at row 112 in the daily price log, shifting the
# existing rows 112-210 down to 113-211, then populate the new row 112 with
# its own values. All other (static) columns mirror the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 112..210 down to 113..211, creating a blank row at 112.
$ws.Rows.Item(112).Insert()

# Fill the newly inserted row 112 with its data.
$ws.Range("A112").Value() = 8
$ws.Range("B112").Value() = "Terminal La Palmera de La Serena"
$ws.Range("C112").Value() = "Coquimbo"
$ws.Range("D112").Value() = 44741
$ws.Range("E112").Value() = 4
$ws.Range("F112").Value() = 100112037
$ws.Range("G112").Value() = "Cebollín"
$ws.Range("H112").Value() = "Sin especificar"
$ws.Range("I112").Value() = "Primera"
$ws.Range("J112").Value() = 1360
$ws.Range("K112").Value() = 1400
$ws.Range("L112").Value() = 1600
$ws.Range("M112").Value() = 1500
$ws.Range("N112").Value() = "$/paquete 6 unidades"
$ws.Range("O112").Value() = "Provincia del Elquí"
$ws.Range("P112").Value() = 250
$ws.Range("Q112").Value() = 6
$ws.Range("R112").Value() = "Hortaliza"
